# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape, and swaps the NEARProtocol/MXToken row order+data, per commit:
# "Updated cryptos list on Sun Nov 19 20:29:38 UTC 2023 with GitHub Actions"
#
# Note: some new Price values (e.g. "245.11", "0.382") would otherwise be
# auto-converted by Excel into numbers; a leading apostrophe forces them to
# stay plain text, matching the original cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "36.887.74"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.979.15"
$ws.Range("E3").Value = "  +0.80%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5: BNB
$ws.Range("D5").Value = "'245.11"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6: XRP
$ws.Range("E6").Value = "  +1.34%  "

# Row 7: Solana
$ws.Range("D7").Value = "'60.82"
$ws.Range("E7").Value = "  +2.88%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.382"
$ws.Range("E9").Value = "  +1.93%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "'0.0799"
$ws.Range("E10").Value = "  -1.97%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +0.72%  "

# Row 12: Chainlink
$ws.Range("D12").Value = "'14.92"
$ws.Range("E12").Value = "  +8.60%  "

# Row 13: Polygon
$ws.Range("E13").Value = "  +1.82%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'22.03"
$ws.Range("E14").Value = "  -0.98%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.273.21"
$ws.Range("E15").Value = "  +0.96%  "

# Row 16: Polkadot
$ws.Range("D16").Value = "'5.48"
$ws.Range("E16").Value = "  +3.98%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "1.977.63"
$ws.Range("E17").Value = "  +0.66%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "36.806.74"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19: Litecoin
$ws.Range("E19").Value = "  +0.34%  "

# Row 20: ShibaInu
$ws.Range("D20").Value = "0.0₃0860"
$ws.Range("E20").Value = "  +0.17%  "

# Row 21: Uniswap
$ws.Range("D21").Value = "'5.17"
$ws.Range("E21").Value = "  +2.06%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'229.84"
$ws.Range("E22").Value = "  +0.28%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.00%  "

# Row 24: PancakeSwap
$ws.Range("E24").Value = "  +2.32%  "

# Row 25: Toncoin
$ws.Range("E25").Value = "  +0.67%  "

# Row 26: Kaspa
$ws.Range("E26").Value = "  +1.97%  "

# Row 27: Cosmos
$ws.Range("E27").Value = "  +0.17%  "

# Row 28: Monero
$ws.Range("D28").Value = "'163.13"
$ws.Range("E28").Value = "  +1.81%  "

# Row 29: EthereumClassic
$ws.Range("E29").Value = "  +0.69%  "

# Row 30: ImmutableX
$ws.Range("E30").Value = "  +17.80%  "

# Row 31: Stellar
$ws.Range("E31").Value = "  +1.55%  "

# Row 32: Filecoin
$ws.Range("E32").Value = "  +3.02%  "

# Row 33: Hedera
$ws.Range("E33").Value = "  +0.29%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").Value = "'4.53"
$ws.Range("E34").Value = "  +5.85%  "

# Row 35: BinanceUSD
$ws.Range("E35").Value = "  -0.01%  "

# Row 36: LidoDAOToken
$ws.Range("D36").Value = "'2.26"
$ws.Range("E36").Value = "  -0.77%  "

# Row 37: RenderToken
$ws.Range("D37").Value = "'3.37"
$ws.Range("E37").Value = "  -0.95%  "

# Row 38: WEMIXToken
$ws.Range("E38").Value = "  +0.28%  "

# Row 39: THORChain
$ws.Range("D39").Value = "'5.54"
$ws.Range("E39").Value = "  -7.80%  "

# Row 40: Cronos
$ws.Range("D40").Value = "'0.0995"
$ws.Range("E40").Value = "  +0.87%  "

# Row 41: HuobiToken
$ws.Range("E41").Value = "  +0.72%  "

# Row 42: TrustWalletToken
$ws.Range("E42").Value = "  +0.58%  "

# Row 43: VeChain
$ws.Range("E43").Value = "  +0.30%  "

# Row 44: InjectiveProtocol
$ws.Range("D44").Value = "'16.42"
$ws.Range("E44").Value = "  +1.25%  "

# Row 45: Maker
$ws.Range("D45").Value = "1.371.86"
$ws.Range("E45").Value = "  +0.80%  "

# Row 46: Aave
$ws.Range("D46").Value = "'90.03"
$ws.Range("E46").Value = "  +2.54%  "

# Row 47: ARBITRUM
$ws.Range("E47").Value = "  -0.16%  "

# Row 48: FraxShare
$ws.Range("D48").Value = "'7.28"
$ws.Range("E48").Value = "  +1.78%  "

# Row 49: NEARProtocol
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.81"
$ws.Range("E49").Value = "  -0.75%  "

# Row 50: MXToken
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.99"
$ws.Range("E50").Value = "  +12.25%  "

# Row 51: MultiversX
$ws.Range("D51").Value = "'46.14"
$ws.Range("E51").Value = "  +5.12%  "
